# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff" on all three
#   sheets (Overview!E2/F2, zh-cn!C2, de-de!C2 all share the same string).
# - Bump the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps forward a bit to reflect the new handoff run.
# - Widen the status columns (they now hold the longer "Ready for handoff"
#   text) to roughly match the post-edit column widths.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- status text ------------------------------------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- timestamps ---------------------------------------------------------
$overview.Range("G2").Value = "2016-08-23 19:02:39"
$dede.Range("H2").Value     = "2016-08-23 19:02:39"
$zhcn.Range("H2").Value     = "2016-08-23 19:02:34"

# --- column widths --------------------------------------------------
# Status columns grew to fit "Ready for handoff". ColumnWidth here is
# expressed in characters (same units Excel's Format > Column Width
# dialog uses); 16.35 is the value that lands the underlying stored
# width on the post-edit figure.
$overview.Columns.Item(5).ColumnWidth = 16.35
$overview.Columns.Item(6).ColumnWidth = 16.35
$zhcn.Columns.Item(3).ColumnWidth     = 16.35
$dede.Columns.Item(3).ColumnWidth     = 16.35
